# BirdColorGenetica.xlsx update
# Extends the color-genetics lookup table with the "Green" and "White"
# body-color cases (Pastel combinations) and updates the selected/active
# cell on the sheet to reflect the newly added search cage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slightly widen columns D and E (Body/Breast results) now that longer
# "Green Pastel" / "White Pastel" values are shown, and drop the old
# "best fit" auto-sizing in favor of the new fixed widths.
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 11.65

# --- Row 26: finish the "Green" combination pair (header half already existed) ---
$ws.Range("B26").Value2 = "Pastel"
$ws.Range("C26").Value2 = "Green"
$ws.Range("D26").Value2 = "Green Pastel"
$ws.Range("E26").Value2 = "Green Pastel"

# --- Row 27: second half of the "Green" combination ---
$ws.Range("B27").Value2 = "Green"
$ws.Range("C27").Value2 = "Pastel"
$ws.Range("D27").Value2 = "Green Pastel"
$ws.Range("E27").Value2 = "Green Pastel"

# --- Row 28: first half of the "White" combination ---
$ws.Range("B28").Value2 = "Pastel"
$ws.Range("C28").Value2 = "White"
$ws.Range("D28").Value2 = "White Pastel"
$ws.Range("E28").Value2 = "White Pastel"

# --- Row 29: second half of the "White" combination ---
$ws.Range("B29").Value2 = "White"
$ws.Range("C29").Value2 = "Pastel"
$ws.Range("D29").Value2 = "White Pastel"
$ws.Range("E29").Value2 = "White Pastel"

# Center-align the newly added rows to match the style used throughout
# the rest of the lookup table.
$ws.Range("B27:E29").HorizontalAlignment = -4108

# Move the active selection to the newly added search cage (bottom-right
# of the table), matching the updated view of the sheet.
$null = $ws.Range("E29").Select()
